# feat: add 2022-Q4 data
#
# 1) Duplicate the existing "2022-Q3" sheet, placing the new copy in front of
#    it (so the tab order becomes ... -> 2022-Q4 -> 2022-Q3 -> ...), then
#    rename the duplicate to "2022-Q4" and refresh it with the new quarter's
#    numbers.
# 2) Prepend a new row on the "总计" (summary) sheet's data table for the
#    2022-Q4 figures, pushing the existing 2022-Q3 / 2021-Q3 rows down by
#    one and renumbering the index column.

$wb = $excel.ActiveWorkbook

# --- 1. New "2022-Q4" worksheet, cloned from "2022-Q3" -----------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Updated fund figures for the new quarter. D2:G2 hold numeric-looking
# text (matching the source sheet's text-typed cells), so a leading
# apostrophe keeps them stored as text instead of being coerced to numbers.
$q4.Range("D2").Value = "'4.89"
$q4.Range("E2").Value = "'99.23"
$q4.Range("F2").Value = "'3.83"
$q4.Range("G2").Value = "'0.1873"
$q4.Range("H2").Value = 8

# --- 2. Prepend a summary row on "总计" for 2022-Q4 ---------------------
$totals = $wb.Worksheets.Item("总计")

# The table is growing from 2 data rows to 3, so A4 is a brand-new cell.
# Clone A3's formatting into it first (matches the bordered/centered index
# column style used by A2/A3), then fill in the shifted rows bottom-up so
# each row's source data is read before it gets overwritten.
$totals.Range("A3").Copy($totals.Range("A4"))

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2021-Q3"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 0.13

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0.18

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.19
